# Adds a new "15. Delete User" test-case section (rows 194-204) to the
# "Test Plan Final" worksheet, mirroring the layout of the existing
# "16. File and data synchronization" section (rows 177-188) but with a
# single test case entry instead of two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan Final")

# ---------------------------------------------------------------------
# 1. Merge the cells for the new block first (so the subsequent format
#    paste lands on the already-merged ranges, matching how the template
#    rows are laid out).
# ---------------------------------------------------------------------
$merges = @(
  "A194:G194",
  "A195:B195", "C195:G195",
  "A196:B196", "C196:G196",
  "A197:B197", "C197:G197",
  "A198:G198",
  "B199:C199", "D199:E199",
  "B200:C200", "D200:E200",
  "A201:B201", "C201:G201",
  "A202:B202", "C202:G202",
  "A203:B203", "C203:G203",
  "A204:B204", "C204:G204"
)
foreach ($m in $merges) {
    $ws.Range($m).MergeCells = $true
}

# ---------------------------------------------------------------------
# 2. Copy the formatting from the equivalent template rows of the
#    "16. File and data synchronization" section onto the new rows.
#    Row 184 (the second test case row of section 16) is skipped because
#    the new section only has a single test case row.
# ---------------------------------------------------------------------
$rowPairs = @(
  @(177,194),
  @(178,195),
  @(179,196),
  @(180,197),
  @(181,198),
  @(182,199),
  @(183,200),
  @(185,201),
  @(186,202),
  @(187,203),
  @(188,204)
)
foreach ($p in $rowPairs) {
    $srcRow = $p[0]
    $dstRow = $p[1]
    $ws.Range("A$srcRow`:G$srcRow").Copy()
    $ws.Range("A$dstRow`:G$dstRow").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 3. Row heights for the rows that use a custom height in the template.
# ---------------------------------------------------------------------
$ws.Rows.Item(194).RowHeight = 20.4
$ws.Rows.Item(195).RowHeight = 15
$ws.Rows.Item(200).RowHeight = 105.6

# ---------------------------------------------------------------------
# 4. Cell values. The brand-new strings are written first, and in the
#    same left-to-right / top-to-bottom order they occur in the row data,
#    so that they are appended to the shared-string table in that order.
# ---------------------------------------------------------------------
$ws.Range("C195").Value = "Test if delete user removes all data related to the user."
$ws.Range("B200").Value = "1. Login to the account.`n2. Go to menu and click delete account`n"
$ws.Range("A194").Value = "15. Delete User"
$ws.Range("A200").Value = "TC15-001"
$ws.Range("D200").Value = "An error message will be printed out on relogin and the user related data are deleted from cloud service (Cognito,S3 and DynamoDB)"
$ws.Range("C203").Value = "2135 hrs"
$ws.Range("F200").Value = "All the user retated data and files are deleted and and error message is printed."

# Remaining cells reuse strings already present in the workbook.
$ws.Range("A195").Value = "Objective"
$ws.Range("A196").Value = "Classification "
$ws.Range("A197").Value = "Pre-requisite (if any)"

$ws.Range("A199").Value = "Case No"
$ws.Range("B199").Value = "Action"
$ws.Range("D199").Value = "Expected result"
$ws.Range("F199").Value = "Actual Result"
$ws.Range("G199").Value = "Pass/Fail/Others"

$ws.Range("G200").Value = "Pass"

$ws.Range("A201").Value = "Tester Name"
$ws.Range("C201").Value = "Abhi Jay Krishnan"

$ws.Range("A202").Value = "Date"
$ws.Range("C202").Value = 43047

$ws.Range("A203").Value = "Time"

$ws.Range("A204").Value = "Pass/Fail/Others"
$ws.Range("C204").Value = "Pass"

# ---------------------------------------------------------------------
# 5. Update the selection to match the author's final cursor position.
# ---------------------------------------------------------------------
$ws.Range("F209").Select()
